# Apply the update described by the commit: append new order rows (32-41)
# to the "Orders" sheet and update the packed "Number" string on "Summary".

$wb = $excel.ActiveWorkbook
$orders = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

function Set-TextCell($range, $value) {
    # Force the cell to be stored as text (matching the rest of the sheet,
    # where numeric-looking values like "16", "10", "1" are kept as text),
    # then drop back to the default "Normal" style so no stray style index
    # (s="...") is left on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# New rows appended after the existing last row (31).
# Columns used: A (PackageID, only on first row of a group), C (FlowerName), F (Number)
$orders.Range("C32").Value = "229_黄蝴蝶_Yellow Butterfly_Rosa rugosa Thunb._10stems"
Set-TextCell $orders.Range("F32") "5"

$orders.Range("C33").Value = "411_紫罗兰白_violet white_undefined_1bunch"
Set-TextCell $orders.Range("F33") "15"

$orders.Range("C34").Value = "412_紫罗兰粉_violet pink_undefined_1bunch"
Set-TextCell $orders.Range("F34") "15"

$orders.Range("C35").Value = "505_紫罗兰紫_violet purple_undefined_1bunch"
Set-TextCell $orders.Range("F35") "10"

Set-TextCell $orders.Range("A36") "7"
$orders.Range("C36").Value = "316_尤加利叶大叶_Eucalyptus Cinerea_undefined_1bunch"
Set-TextCell $orders.Range("F36") "20"

$orders.Range("C37").Value = "70_朝霞mini_undefined_Gerbera L._20stems"
Set-TextCell $orders.Range("F37") "15"

$orders.Range("C38").Value = "71_霜雪mini_Snowy_Gerbera L._20stems"
Set-TextCell $orders.Range("F38") "10"

$orders.Range("C39").Value = "463_玉兰枝_magnolia flower`nwhite/purple_undefined_1bunch"
Set-TextCell $orders.Range("F39") "10"

$orders.Range("C40").Value = "328_卢荀草_undefined_undefined_1bunch"
Set-TextCell $orders.Range("F40") "10"

$orders.Range("C41").Value = "77_珍爱mini_undefined_Gerbera L._20stems"

# Update the Summary sheet's packed Number string (G2) to include the
# newly-added rows' Number values appended to the existing string.
Set-TextCell $summary.Range("G2") "016111210101614101030101010305040105201010101515102010355515151020151010100"
